$d = $word.ActiveDocument

# Locate the paragraph that contains "Ver no Jupiter Salvar em pdf Salvar em docx".
# The edit removes that paragraph, the "(c) 2020 ... Creative Commons
# Attribution" paragraph right after it, and the single blank paragraph
# that precedes both of them (directly after the "MacGrall-Hill"
# bibliography line). The bibliography paragraph itself, and the blank
# paragraph that originally followed the removed block, are left intact.

$target = "Ver no Jupiter Salvar em pdf Salvar em docx"
$footer = "Creative Commons Attribution"

$count = $d.Paragraphs.Count
$targetIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($target)) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Paragraph right before the "Ver no Jupiter..." line is the blank
    # separator paragraph that should also be removed.
    $blankIndex = $targetIndex - 1

    # Paragraph right after should contain the copyright/footer text;
    # delete from the last one down to the first so indices of the
    # paragraphs being removed stay valid as each delete happens.
    $footerIndex = $targetIndex + 1
    if (-not $d.Paragraphs.Item($footerIndex).Range.Text.Contains($footer)) {
        $footerIndex = $targetIndex
    }

    if ($footerIndex -ne $targetIndex) {
        $d.Paragraphs.Item($footerIndex).Range.Delete()
    }
    $d.Paragraphs.Item($targetIndex).Range.Delete()
    if ($blankIndex -ge 1 -and $d.Paragraphs.Item($blankIndex).Range.Text.Trim().Length -eq 0) {
        $d.Paragraphs.Item($blankIndex).Range.Delete()
    }
}
